$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking cells, then set values, so they are stored as text
# (matching the original inline-string / text representation of these cells).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "241.86"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "21"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.91"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "21"
$ws.Range("B4").Value = "HuobiToken"
$ws.Range("C4").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.381"
$ws.Range("E4").Value = "3HuobiTokenHT"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "21"
$ws.Range("B5").Value = "Cronos"
$ws.Range("C5").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05944"
$ws.Range("E5").Value = "4CronosCRO"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "21"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.395"
$ws.Range("E6").Value = "5GateTokenGT"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "21"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.459"
$ws.Range("E7").Value = "6KuCoinTokenKCS"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "21"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8040"
$ws.Range("E8").Value = "7MXTokenMX"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "21"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9119"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "21"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1411"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "21"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07415"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "21"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03307"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "21"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03033"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "21"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09329"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "21"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.874"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "21"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001574"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "21"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04525"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "21"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005944"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "21"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006129"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "21"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004995"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "21"
$ws.Range("B21").Value = "UpBots"
$ws.Range("C21").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.007494"
$ws.Range("E21").Value = "20UpBotsUBXTBestin24h"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "21"
$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0009859"
$ws.Range("E22").Value = "21BitKanKAN"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "21"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.00007803"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "21"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.614"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "21"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "21"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "21"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "21"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "21"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "21"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "21"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "21"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "21"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "21"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "21"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "21"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "21"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "21"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "21"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "21"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03850"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "21"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006068"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "21"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1063"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "21"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002801"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "21"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007201"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "21"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005192"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "21"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "21"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005804"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "21"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "21"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002261"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "21"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "21"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "21"
